# Insert a new column at C ("statut_name") and shift existing columns C:L -> D:M.
# Resulting layout:
#   A statut | B statut_label | C statut_name | D NCTId | E eudraCT | F CTIS
#   G completion_year | H clinical_trial_title | I acronym | J results_1y
#   K results_3y | L results | M intervention_type

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C; this shifts the old C:L block to D:M.
$ws.Columns.Item(3).Insert()

# Copy the header formatting from the (now shifted) NCTId header in D1 onto the
# new statut_name header in C1, so it keeps the bold/centered/bordered look.
$ws.Cells.Item(1, 4).Copy()
$ws.Cells.Item(1, 3).PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Set the new header text.
$ws.Cells.Item(1, 3).Value = "statut_name"

# Determine the last used data row (header is row 1, data starts row 2).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

# Fill the new column's data rows with the fixed status label.
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 3).Value = "pas de résultat ni de publication"
}
